$wb = $excel.ActiveWorkbook

# --- Sheet "OM" ---
$ws = $wb.Worksheets.Item("OM")
$ws.Range("B2").Value = 67
$ws.Range("B3").Value = 28.3
$ws.Range("B4").Value = 4.3

# --- Sheet "NV" ---
$ws = $wb.Worksheets.Item("NV")
$ws.Range("B2").Value = 89
$ws.Range("B3").Value = 9.9
$ws.Range("B4").Value = 1
# New row 5: A5 must stay a text label ("3"), like the other A-column labels.
$ws.Range("A5").Value = "'3"
$ws.Range("B5").Value = 0.1

# --- Sheet "NR" ---
$ws = $wb.Worksheets.Item("NR")
$ws.Range("B2").Value = 35.6
$ws.Range("B3").Value = 21.9
$ws.Range("B5").Value = 10
$ws.Range("B6").Value = 6.7
$ws.Range("B7").Value = 4.9
$ws.Range("B8").Value = 2.2
$ws.Range("B9").Value = 1.3
$ws.Range("B10").Value = 1.4
$ws.Range("B11").Value = 0.4
$ws.Range("B12").Value = 0.5
$ws.Range("B14").Value = 0.2
$ws.Range("B15").Value = 0.2

# --- Sheet "ALL" ---
$ws = $wb.Worksheets.Item("ALL")
$ws.Range("B2").Value = 20.2
$ws.Range("B3").Value = 23.5
$ws.Range("B4").Value = 19.9
$ws.Range("B5").Value = 13.1
$ws.Range("B6").Value = 7.5
$ws.Range("B7").Value = 6.6
$ws.Range("B8").Value = 3.2
$ws.Range("B9").Value = 2.2
$ws.Range("B10").Value = 1.6
$ws.Range("B11").Value = 0.5
$ws.Range("B13").Value = 0.5
$ws.Range("B14").Value = 0.4
$ws.Range("B15").Value = 0.2

# --- Sheet "summary" ---
$ws = $wb.Worksheets.Item("summary")
$ws.Range("C2").Value = 1.86
$ws.Range("D2").Value = 2.36
$ws.Range("C3").Value = 2.27
$ws.Range("D3").Value = 2.29
# B6 switches from the "2" label to the "3" label (text, like its neighbours).
$ws.Range("B6").Value = "'3"
